# Apply cryptos list price/volume refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.427.77'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '2.104.88'
$ws.Range('E3').Value = '  +1.09%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').Value = '''334.40'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').Value = '''1.003'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').Value = '''0.5226'
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('D8').Value = '''0.4552'
$ws.Range('E8').Value = '  +5.56%  '
$ws.Range('D9').Value = '''53.23'
$ws.Range('E9').Value = '  +15.54%  '
$ws.Range('D10').Value = '''0.08934'
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('D11').Value = '''1.179'
$ws.Range('E11').Value = '  +1.64%  '
$ws.Range('D12').Value = '''24.24'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').Value = '2.101.21'
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').Value = '''6.850'
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('E15').Value = '  +4.95%  '
$ws.Range('D16').Value = '''96.55'
$ws.Range('E16').Value = '  +1.39%  '
$ws.Range('D17').Value = '''0.00001142'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '''1.006'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '''0.06649'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '''19.23'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').Value = '''6.337'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '30.486.84'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').Value = '''12.50'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').Value = '''2.360'
$ws.Range('E25').Value = '  +3.06%  '
$ws.Range('D26').Value = '2.363.52'
$ws.Range('E26').Value = '  +1.72%  '
$ws.Range('D27').Value = '''22.27'
$ws.Range('D28').Value = '''2.541'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('D29').Value = '''162.69'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('D30').Value = '''133.13'
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D31').Value = '''1.214'
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('D32').Value = '''0.1073'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').Value = '''1.668'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('D34').Value = '''6.374'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('D35').Value = '''3.939'
$ws.Range('E35').Value = '  +3.03%  '
$ws.Range('E36').Value = '  +6.60%  '
$ws.Range('D37').Value = '''5.768'
$ws.Range('E37').Value = '  +5.90%  '
$ws.Range('D38').Value = '''0.02592'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').Value = '''0.06846'
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('D40').Value = '''0.2303'
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('D41').Value = '''12.72'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D43').Value = '''1.253'
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''14.13'
$ws.Range('E44').Value = '  +1.47%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '''2.321'
$ws.Range('E45').Value = '  +5.46%  '
$ws.Range('D46').Value = '''0.6380'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('D47').Value = '''3.662'
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').Value = '''1.251'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').Value = '''0.00000000349'
$ws.Range('E49').Value = '  +21.00%  '
$ws.Range('D50').Value = '''0.3418'
$ws.Range('E50').Value = '  +24.56%  '
$ws.Range('D51').Value = '''83.49'
$ws.Range('E51').Value = '  +2.30%  '
